$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.933.69'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.816.97'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4691'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.74%  '
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07375'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8730'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.40'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = '1.814.94'
$ws.Range('E12').Value = '  +6.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.386'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07085'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.94'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '26.972.06'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.334'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').Value = '2.040.23'
$ws.Range('E24').Value = '  +4.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.890'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '150.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.182'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08958'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7679'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.168'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.516'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.910'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.087'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01965'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05297'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.970'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.282'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5355'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.332'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1656'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.467'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4938'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('E51').Value = '  -0.14%  '
